$wb = $excel.ActiveWorkbook

# ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1132.6666
$ws.Range("I6").Value = 700
$ws.Range("J6").Value = 1219.2
$ws.Range("K6").Value = 2100
$ws.Range("L6").Value = 3657.6
$ws.Range("M6").Value = -1988
$ws.Range("N6").Value = -3881.6

# ALC row 8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 281.3387
$ws.Range("I8").Value = 42.5
$ws.Range("J8").Value = 297.81033
$ws.Range("K8").Value = 127.5
$ws.Range("L8").Value = 893.4309900000001
$ws.Range("M8").Value = 11.5
$ws.Range("N8").Value = -1171.43099

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 840.2917
$ws.Range("I33").Value = 831.4211
$ws.Range("J33").Value = 874
$ws.Range("K33").Value = 831.4211
$ws.Range("L33").Value = 874
$ws.Range("M33").Value = -602.4211
$ws.Range("N33").Value = -1332

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 29690.047
$ws.Range("I135").Value = 1871.3334
$ws.Range("K135").Value = 16842.0006
$ws.Range("M135").Value = -14307.0006

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 30119.6
$ws.Range("I5").Value = 30119.6
$ws.Range("K5").Value = 30119.6
$ws.Range("M5").Value = -30007.6

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27936.717
$ws.Range("I32").Value = 29977.275
$ws.Range("K32").Value = 29977.275
$ws.Range("M32").Value = -29690.275

# ARM row 37
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 23124.625
$ws.Range("J37").Value = 31249.25
$ws.Range("L37").Value = 31249.25
$ws.Range("N37").Value = -31795.25

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2746.7334
$ws.Range("I61").Value = 1304.12
$ws.Range("K61").Value = 1304.12
$ws.Range("M61").Value = -1092.12

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 502044.9
$ws.Range("I74").Value = 601203.9
$ws.Range("K74").Value = 601203.9
$ws.Range("M74").Value = -600329.9

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 502044.9
$ws.Range("I77").Value = 601203.9
$ws.Range("K77").Value = 3006019.5
$ws.Range("M77").Value = -3001651.5

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2065
$ws.Range("I122").Value = 1910.5555
$ws.Range("K122").Value = 5731.666499999999
$ws.Range("M122").Value = -3281.666499999999

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1147.4615
$ws.Range("I132").Value = 1075.25
$ws.Range("J132").Value = 2014
$ws.Range("K132").Value = 3225.75
$ws.Range("L132").Value = 6042
$ws.Range("M132").Value = -695.75
$ws.Range("N132").Value = -11102

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2746.7334
$ws.Range("I136").Value = 1304.12
$ws.Range("K136").Value = 3912.36
$ws.Range("M136").Value = -1362.36

# ARM row 138
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 30119.6
$ws.Range("I4").Value = 30119.6
$ws.Range("K4").Value = 30119.6
$ws.Range("M4").Value = -30004.6

# BSM row 40
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 63999
$ws.Range("J40").Value = 63999
$ws.Range("L40").Value = 63999
$ws.Range("N40").Value = -64529

# BSM row 96
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 39499.5
$ws.Range("J96").Value = 63999
$ws.Range("L96").Value = 63999
$ws.Range("N96").Value = -69491

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 100264.13
$ws.Range("J140").Value = 100264.13
$ws.Range("L140").Value = 100264.13
$ws.Range("N140").Value = -110624.13

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 306
$ws.Range("I22").Value = 268
$ws.Range("J22").Value = 401
$ws.Range("K22").Value = 268
$ws.Range("L22").Value = 401
$ws.Range("M22").Value = 82
$ws.Range("N22").Value = -1101

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1069.2
$ws.Range("I122").Value = 910.2222
$ws.Range("K122").Value = 2730.6666
$ws.Range("M122").Value = -280.6666

# CRP row 133
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 164229.5
$ws.Range("J133").Value = 208163
$ws.Range("L133").Value = 208163
$ws.Range("N133").Value = -213223

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1816.8
$ws.Range("J5").Value = 994.6667
$ws.Range("L5").Value = 2984.0001
$ws.Range("N5").Value = -3208.0001

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 259.14285
$ws.Range("I17").Value = 232.4
$ws.Range("J17").Value = 326
$ws.Range("K17").Value = 697.2
$ws.Range("L17").Value = 978
$ws.Range("M17").Value = -528.2
$ws.Range("N17").Value = -1316

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4210.923
$ws.Range("J68").Value = 4299.36
$ws.Range("L68").Value = 12898.08
$ws.Range("N68").Value = -14520.08

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 4210.923
$ws.Range("J71").Value = 4299.36
$ws.Range("L71").Value = 38694.24
$ws.Range("N71").Value = -46806.24

# CUL row 74
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# CUL row 77
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1816.8
$ws.Range("J135").Value = 994.6667
$ws.Range("L135").Value = 8952.0003
$ws.Range("N135").Value = -14022.0003

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2143.2104
$ws.Range("I137").Value = 2078.2942
$ws.Range("J137").Value = 2695
$ws.Range("K137").Value = 6234.882599999999
$ws.Range("L137").Value = 8085
$ws.Range("M137").Value = -1134.882599999999
$ws.Range("N137").Value = -18285

# CUL row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 12200.75
$ws.Range("J138").Value = 15591
$ws.Range("L138").Value = 46773
$ws.Range("N138").Value = -57053

# GSM row 40
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 19117.25
$ws.Range("I40").Value = 17500
$ws.Range("J40").Value = 19656.334
$ws.Range("K40").Value = 17500
$ws.Range("L40").Value = 19656.334
$ws.Range("M40").Value = -17349
$ws.Range("N40").Value = -19958.334

# GSM row 62
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# GSM row 65
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# GSM row 141
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 86700
$ws.Range("J141").Value = 86700
$ws.Range("L141").Value = 86700
$ws.Range("N141").Value = -97060

# LTW row 13
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 26665
$ws.Range("I13").Value = 32997.5
$ws.Range("J13").Value = 14000
$ws.Range("K13").Value = 32997.5
$ws.Range("L13").Value = 14000
$ws.Range("M13").Value = -32857.5
$ws.Range("N13").Value = -14280

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6149.0557
$ws.Range("I46").Value = 2962.875
$ws.Range("J46").Value = 8698
$ws.Range("K46").Value = 2962.875
$ws.Range("L46").Value = 8698
$ws.Range("M46").Value = -2774.875
$ws.Range("N46").Value = -9074

# LTW row 74
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 46833.168
$ws.Range("J74").Value = 55249.75
$ws.Range("L74").Value = 55249.75
$ws.Range("N74").Value = -57245.75

# LTW row 77
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 46833.168
$ws.Range("J77").Value = 55249.75
$ws.Range("L77").Value = 165749.25
$ws.Range("N77").Value = -175733.25

# LTW row 121
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 40000
$ws.Range("J121").Value = 40000
$ws.Range("L121").Value = 40000
$ws.Range("N121").Value = -43494

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 17935
$ws.Range("I136").Value = 19902.412
$ws.Range("K136").Value = 59707.236
$ws.Range("M136").Value = -57157.236
